$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 98.84999999999999
$ws.Range("H2").Value = 0.9399999999999999

# Row 4
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 98.45999999999999
$ws.Range("H4").Value = 1.44

# Row 5
$ws.Range("C5").Value = 99.62
$ws.Range("D5").Value = 0.38
$ws.Range("E5").Value = 0.03
$ws.Range("F5").Value = 99.62
$ws.Range("G5").Value = 99.23
$ws.Range("H5").Value = 0.9399999999999999
